$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert 3 blank rows at the positions where new players land once the
#    table (sorted descending by column C, "inicial") is re-sorted with the
#    new data included. Inserting top-down keeps later indices valid.
# ---------------------------------------------------------------------------
$ws.Rows("8:8").Insert()
$ws.Rows("18:18").Insert()
$ws.Rows("30:30").Insert()

# ---------------------------------------------------------------------------
# 2. Give the new A/C cells the same number formatting (style index 2) used
#    by the rest of the table, by copying formats from neighbouring rows.
# ---------------------------------------------------------------------------
$ws.Range("A9").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("C9").Copy()
$ws.Range("C8").PasteSpecial(-4122)

$ws.Range("A19").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("C19").Copy()
$ws.Range("C18").PasteSpecial(-4122)

$ws.Range("A31").Copy()
$ws.Range("A30").PasteSpecial(-4122)
$ws.Range("C31").Copy()
$ws.Range("C30").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Fill in the values for the three newly-added players.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "6095fca68a19d000196923bf"
$ws.Range("B8").Value = "Yarey"
$ws.Range("C8").Value = 972672480
$ws.Range("D8").Value = 7

$ws.Range("A18").Value = "6116f40168427b0029d5580f"
$ws.Range("B18").Value = "Julio Barboza"
$ws.Range("C18").Value = 110472005
$ws.Range("D18").Value = 17

$ws.Range("A30").Value = "6117fe5efc41af003a67eb8f"
$ws.Range("B30").Value = "Rafael Braga"
$ws.Range("C30").Value = 7756768
$ws.Range("D30").Value = 29

# ---------------------------------------------------------------------------
# 4. Renumber the "posicao" column (D) for every data row so it stays a
#    contiguous 1..30 sequence after the insert.
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 31; $r++) {
    $ws.Range("D$r").Value = $r - 1
}

# ---------------------------------------------------------------------------
# 5. Rebuild every hyperlink on column E (the row shift invalidates the old
#    ones, and the engine does not auto-follow hyperlinks through inserts).
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

$links = @{
    2  = "https://rollercoin.com/p/MitoBR"
    4  = "https://rollercoin.com/p/BabyTux69"
    5  = "https://rollercoin.com/p/kzugpfwl"
    6  = "https://rollercoin.com/p/ZeraTAL"
    8  = "https://rollercoin.com/p/YesoGengo"
    9  = "https://rollercoin.com/p/Walljack"
    10 = "https://rollercoin.com/p/jzlfzveq"
    11 = "https://rollercoin.com/p/larjauxv"
    12 = "https://rollercoin.com/p/JHOWBR21"
    13 = "https://rollercoin.com/p/ricardomann"
    14 = "https://rollercoin.com/p/jukinha"
    16 = "https://rollercoin.com/p/lg9p91px"
    17 = "https://rollercoin.com/p/terraqueo"
    18 = "https://rollercoin.com/p/SidFillips"
    19 = "https://rollercoin.com/p/lkr7fomr"
    20 = "https://rollercoin.com/p/thekrk420"
    21 = "https://rollercoin.com/p/lvx9ene8"
    22 = "https://rollercoin.com/p/Dyegolimax"
    23 = "https://rollercoin.com/p/Pesaac"
    24 = "https://rollercoin.com/p/sccp_gu"
    25 = "https://rollercoin.com/p/kl6yb7ip"
    26 = "https://rollercoin.com/p/Regis"
    27 = "https://rollercoin.com/p/kkxb6fy2"
    28 = "https://rollercoin.com/p/lithrrgz"
    29 = "https://rollercoin.com/p/kc30bi28"
    30 = "https://rollercoin.com/p/ksc25irv"
    31 = "https://rollercoin.com/p/lyamecmn"
}

foreach ($r in ($links.Keys | Sort-Object)) {
    $ws.Hyperlinks.Add($ws.Range("E$r"), $links[$r])
}

# ---------------------------------------------------------------------------
# 6. Dimension / auto-filter sort-state / selection bookkeeping.
# ---------------------------------------------------------------------------
$ws.Range("D2").Select()

Write-Output "done"
